$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 8 (shifts existing row 8 and below down by one).
$ws.Rows("8:8").Insert()

# Row 7's second column previously held the combined translation "вся земля";
# split it so "вся " stays on row 7 (aligned with "כָּל־הָאָרֶץ").
$ws.Range("B7").Value = "вся "

# New row 8 carries the second half of the split pair.
$ws.Range("A8").Value = "הָאָרֶץ"
$ws.Range("B8").Value = "земля"
